$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shorten the long "Inscription" (program/orientation) labels to their
# short forms. Using Find/Replace across the used range keeps every row
# that shares the original string pointing at the same (now-renamed)
# shared string, instead of spawning a brand new one per edited cell.
$ws.Cells.Replace("Bachelier en informatique et systèmes orientation, réseaux et télécommunicationsBloc 1", "Informatique")
$ws.Cells.Replace("Bachelier en techniques graphiques, orientation techniques infographiques", "Infographie")
$ws.Cells.Replace("Bachelier en électronique, orientation électronique appliquée", "Electronique")
$ws.Cells.Replace("Bachelier en biotechnique(bloc1)", "Biotechnique")

# Update the sheet's remembered selection/active cell.
$ws.Range("F13").Select()
